$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts existing rows 2-52 down to 3-53)
$ws.Rows.Item(2).Insert()

# Copy date-column formatting from the row below into the new A2 cell
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B2:D2").ClearFormats()
$ws.Range("E2").Clear()

# New first data row values
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 7.226520411029069
$ws.Range("D2").Value = 2008
$ws.Range("C4").Value = 4.268860212333636
$ws.Range("C6").Value = -7.266312015249776
$ws.Range("C7").Value = 3.184002331674129
$ws.Range("E7").Value = 6.704254199558113
$ws.Range("C8").Value = 6.958243460951929
$ws.Range("E8").Value = 12.21658306395068
$ws.Range("C9").Value = 8.626810748872327
$ws.Range("E9").Value = 4.739201070534826
$ws.Range("C10").Value = 9.469137444079934
$ws.Range("E10").Value = 8.079264579851909
$ws.Range("C11").Value = 3.449685446853534
$ws.Range("E11").Value = 3.26507595662513
$ws.Range("C12").Value = 3.358206407534947
$ws.Range("E12").Value = 4.701432377325987
$ws.Range("C13").Value = -1.480934717826909
$ws.Range("E13").Value = 0.7772706050320544
$ws.Range("C14").Value = 0.3081076735359067
$ws.Range("E14").Value = 3.972902167062387
$ws.Range("C15").Value = 5.427992542801308
$ws.Range("E15").Value = 4.945882057432871
$ws.Range("C16").Value = 3.901355411819707
$ws.Range("E16").Value = 4.658857392675264
$ws.Range("C17").Value = 4.970284184513551
$ws.Range("E17").Value = 4.488174889976171
$ws.Range("C18").Value = 5.331683351557981
$ws.Range("E18").Value = 4.089819750351786
$ws.Range("C19").Value = 4.039484738713828
$ws.Range("E19").Value = 4.214976960249173
$ws.Range("C20").Value = 3.254758369308375
$ws.Range("E20").Value = 2.313009565865753
$ws.Range("C21").Value = 4.589070866863865
$ws.Range("E21").Value = 3.829046580278361
$ws.Range("C22").Value = 5.246209615995667
$ws.Range("E22").Value = 4.784022165496182
$ws.Range("C23").Value = 6.011890504679696
$ws.Range("E23").Value = 4.234360353587641
$ws.Range("C24").Value = 3.625873842174787
$ws.Range("E24").Value = 2.330842103296149
$ws.Range("C25").Value = 4.899902276557011
$ws.Range("E25").Value = 3.75342745845737
$ws.Range("C26").Value = 4.86255966374296
$ws.Range("E26").Value = 4.112897401876747
$ws.Range("C27").Value = 3.660106318836931
$ws.Range("E27").Value = 3.270208315717005
$ws.Range("C28").Value = 3.500574054404404
$ws.Range("E28").Value = 3.21661481720994
$ws.Range("C29").Value = 2.983312281417039
$ws.Range("E29").Value = 2.428295356218069
$ws.Range("C30").Value = 2.764740011159428
$ws.Range("E30").Value = 1.643374185611401
$ws.Range("C31").Value = 2.096953540210977
$ws.Range("E31").Value = 3.169670668618951
$ws.Range("C32").Value = -0.9913189363815245
$ws.Range("E32").Value = 1.183532150252908
$ws.Range("C33").Value = -7.578477024949737
$ws.Range("E33").Value = -5.743787238149123
$ws.Range("C34").Value = -7.260793671746435
$ws.Range("E34").Value = 0.00562230452727519
$ws.Range("C35").Value = 0.4989366167094333
$ws.Range("E35").Value = 2.690694906265412
$ws.Range("C36").Value = 4.507091823899212
$ws.Range("E36").Value = 5.429743376942153
$ws.Range("C37").Value = 4.379227219808146
$ws.Range("E37").Value = 4.954652839642848
$ws.Range("C38").Value = 4.097586525396268
$ws.Range("E38").Value = 3.9116372951149
$ws.Range("C39").Value = 7.041577295022128
$ws.Range("E39").Value = 3.388682041315016
$ws.Range("C40").Value = 8.053468068361846
$ws.Range("E40").Value = 3.974997080343634
$ws.Range("C41").Value = 7.397318165265498
$ws.Range("E41").Value = 3.367096865515662
$ws.Range("C42").Value = 7.824284864703746
$ws.Range("E42").Value = 2.586378346096296
$ws.Range("C43").Value = 0.9995490351194292
$ws.Range("E43").Value = 2.834404338648921
$ws.Range("C44").Value = 0.2714278794373248
$ws.Range("E44").Value = 2.460471645027118
$ws.Range("C45").Value = -0.3046246622258053
$ws.Range("E45").Value = 1.976476469605681
$ws.Range("C46").Value = -1.24502235313334
$ws.Range("E46").Value = -1.561801765212567
$ws.Range("C47").Value = -2.798317913999848
$ws.Range("E47").Value = 2.363509743917169
$ws.Range("C48").Value = -2.107534670984712
$ws.Range("E48").Value = 2.747596279389564
$ws.Range("C49").Value = -2.567041707495976
$ws.Range("E49").Value = 1.835066812373642
$ws.Range("C50").Value = -1.735114423676209
$ws.Range("E50").Value = 2.409056355286521
$ws.Range("C51").Value = 2.450219408996213
$ws.Range("E51").Value = 2.677741483899121
$ws.Range("C52").Value = 1.552685227480533
$ws.Range("E52").Value = 2.496145622272206
$ws.Range("C53").Value = 2.64031107104763
$ws.Range("E53").Value = 3.123685491361705
